$d = $word.ActiveDocument

# 1) Ativação date: 2018 -> 2025
$d.Content.Find.Execute(
    "Ativação: 01/01/2018", $false, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# 2) "Programa resumido" section — two paragraphs with identical original
#    text but different replacements (plain PT paragraph, then italic EN-ish
#    paragraph that becomes the English translation). Target each paragraph
#    individually by its Find range to avoid touching the wrong one.
$p11 = $d.Paragraphs.Item(11).Range
$p11.Find.Execute(
    "Introdução a física, Cinemática, Dinâmica, Trabalho, Torque e Momento Angular.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Introdução a física, Cinemática, Dinâmica, Energia, Momento linear, Rotação.", 2) | Out-Null

$p12 = $d.Paragraphs.Item(12).Range
$p12.Find.Execute(
    "Introdução a física, Cinemática, Dinâmica, Trabalho, Torque e Momento Angular.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Introduction to Physics, Kinematics, Dynamics, Energy, Linear momentum, Rotation", 2) | Out-Null

# 3) "Programa" section — detailed Portuguese program text
$d.Content.Find.Execute(
    "1) Introdução a Física: noções de algarismos, análise dimensional, sistemas de unidades.2) Cinemática: leis de Newton e aplicações.3) Trabalho: conservação de energia, forças conservativas, aplicações.4) Impulso: momento linear e conservação de momento linear.5) Torque e momento angular: conservação de momento angular, pêndulo.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "1) Introdução a Física: sistemas de unidades, revisão de vetores, análise dimensional.2) Cinemática: movimento unidimensional, queda livre, movimento bidimensional, projéteis. 3) Dinâmica: leis de Newton, forças, força de atrito, força de resistência do ar, velocidade terminal, movimento circular uniforme, gravitação, aplicações.4) Energia: trabalho, forças conservativas, conservação de energia mecânica, atrito, aplicações.5)  Momento linear: centro de massa, sistema de partículas, conservação do momento linear, colisões, impulso.6) Rotação: variáveis do movimento rotacional, energia cinética rotacional, momento de inércia, torque, rolamento, conservação do momento angular.",
    2) | Out-Null

# 4) "Programa" section — detailed English program text
$d.Content.Find.Execute(
    "1) Introduction to Physics:  significant algharisms, dimensional analysis, units systems.2) Kinematics: Newton's laws and applications.3) Work: energy conservation, conservative forces, applications.4) Impulse: momentum and conservation.5) Torque and Angular Momentum: angular momentum conservation, pendulum.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "1) Introduction to Physics: unit systems, review of vectors, dimensional analysis. 2) Kinematics: one dimensional motion, free fall, bidimensional motion, projectile.  3) Dynamics: Newton’s laws, friction force, drag force, terminal speed, uniform circular motion, gravitation, applications.4) Energy: work, conservative forces, mechanical energy conservation, friction, applications.5)  Linear momentum: center of mass, system of particles, conservation of linear momentum, collisions, impulse.6) Rotation: rotational variables, kinetic energy of rotation, rotational inertia, torque, rolling, conservation of angular momentum",
    2) | Out-Null

# 5) Método (grading method) text
$d.Content.Find.Execute(
    "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "As avaliações serão compostas por provas, projetos, seminários e outras formas que serão utilizadas para a composição das notas. A média final (NF) é calculada pela média simples das notas (N), levando em conta o número n de avaliações, sendo no mínimo duas avaliações: NF= (N1+...+Nn)/n.",
    2) | Out-Null

# 6) Critério (passing criterion) text
$d.Content.Find.Execute(
    "NF≥ 5,0.", $false, $false, $false, $false, $false,
    $true, 1, $false, "NF ≥ 5,0", 2) | Out-Null

# 7) Norma de recuperação (makeup exam rule) text
$d.Content.Find.Execute(
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "(NF+REC)/2 ≥ 5,0, onde REC é uma prova de recuperação a ser aplicada, seguindo as regras da EEL.",
    2) | Out-Null

# 8) Bibliografia — reordered reference list
$d.Content.Find.Execute(
    "NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.1, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008).",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "HALLIDAY, D; RESNICK, R. Fundamentos de Física. Vol.1, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008).NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).",
    2) | Out-Null

Write-Output "done"
